# Generate Report for Archive
# Update localization status for files that have moved from "Ready for handoff"
# into "In Translation": 25b180ed-7e66-4921-838d-c1ea32aa5b85.md (row 3) and
# 719e660d-6d5d-476c-8e53-25d9bdd67da3.md (row 4). The third handed-off file,
# 73baee7e-d827-4a80-9bc6-a83dc0ac0602.md (row 5), stays "Ready for handoff".

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"

# Per-locale detail sheets: Status lives in column C (table column "Status").
foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("C3").Value = $newStatus
    $ws.Range("C4").Value = $newStatus
}

# Overview sheet: the zh-cn / de-de columns (B and C) mirror each locale's status.
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = $newStatus
$overview.Range("C3").Value = $newStatus
$overview.Range("B4").Value = $newStatus
$overview.Range("C4").Value = $newStatus
